$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly (Week 9) activity rows appended below the existing log, rows 355-363.
$data = @(
    ,("Phil",     45513, "Run",     18, 2.25,  125, 0,  2,  5, 9, 0, "Sauntering Hippo", 9)
    ,("Phil",     45513, "Workout", 58, 0,      0, 35, 22,  0, 0, 0, "Sauntering Hippo", 9)
    ,("Matt",     45513, "Ride",    42, 12.67,  0, 13, 30,  0, 0, 0, "Agile Antelope",   9)
    ,("Jeremiah", 45514, "Workout", 33, 0,      0, 11, 16,  6, 0, 0, "Agile Antelope",   9)
    ,("Matt",     45514, "Run",     50, 4.76,  203, 1, 19, 17, 9, 1, "Agile Antelope",   9)
    ,("Eric",     45514, "Run",     12, 1.26,   46, 0,  7,  0, 0, 0, "Wily Hyena",       9)
    ,("Eric",     45514, "Workout", 48, 0,      0, 13, 13, 18, 4, 0, "Wily Hyena",       9)
    ,("Matt",     45514, "Walk",    23, 1.03,   85, 16,  0,  0, 0, 0, "Agile Antelope",   9)
    ,("Steven",   45514, "Walk",    36, 1.8,    79, 36,  0,  0, 0, 0, "Brave Leopard",    9)
)

$startRow = 355
for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $startRow + $r
    $vals = $data[$r]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
    $ws.Cells.Item($row, 9).Value = $vals[8]
    $ws.Cells.Item($row, 10).Value = $vals[9]
    $ws.Cells.Item($row, 11).Value = $vals[10]
    $ws.Cells.Item($row, 12).Value = $vals[11]
    $ws.Cells.Item($row, 13).Value = $vals[12]
}

# Match the date formatting (style) already used by column B ("m/d/yyyy") for the
# newly added rows, by copying the format from the last pre-existing date cell.
$ws.Range("B354").Copy()
$ws.Range("B355:B363").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the new selection / active cell at the bottom of the sheet.
$ws.Range("A364").Select()
